$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Ballot($row, $voter, $cols, $nvotes, $source, $date) {
    $ws.Range("A$row").Value = $voter
    foreach ($c in $cols) {
        $ws.Range("$c$row").Value = "x"
    }
    $ws.Range("AK$row").Value = $nvotes
    $ws.Range("AL$row").Value = $source
    $ws.Range("AM$row").Value = $date
    # Copy the date number format from an existing dated cell so we reuse
    # the workbook's existing style (numFmtId 14) instead of creating a new one.
    $ws.Range("AM45").Copy()
    $ws.Range("AM$row").PasteSpecial(-4122)
}

# New ballots added 12/15-12/16 (commit: "add 5 new ballots from 12/15-12/16")
Set-Ballot 46 "Mark Newman"     @("C","D","E","I","K","M","O","P","Q","V")     10 "Twitter" 43450
Set-Ballot 47 "Richard Justice" @("C","D","E","I","K","O","P","Q","U","V")     10 "Email"   43449
Set-Ballot 48 "Mark Hale"       @("C","D","E","I","K","N","O","Q","V")         9  "Twitter" 43449
Set-Ballot 49 "Jeff Blair"      @("C","D","E","O","V")                         5  "Twitter" 43449

# Update the frozen-pane view so the active selection matches the new last row.
$win = $excel.ActiveWindow
[void]($win.FreezePanes = $false)
[void]($ws.Range("B2").Select())
[void]($win.FreezePanes = $true)
[void]($ws.Range("B49").Select())
